$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Range("G3").Value = 1.41
$ws.Range("I3").Value = 14.5

# Row 4
$ws.Range("I4").Value = 2.24
$ws.Range("J4").Value = 3.25
$ws.Range("K4").Value = 3.45
$ws.Range("P4").Value = 1.66

# Row 5
$ws.Range("F5").Value = 2.12
$ws.Range("G5").Value = 2.32
$ws.Range("H5").Value = 3.65
$ws.Range("I5").Value = 4.3
$ws.Range("J5").Value = 3.4
$ws.Range("K5").Value = 3.7

# Row 6
$ws.Range("H6").Value = 3.35
$ws.Range("I6").Value = 3.8
$ws.Range("K6").Value = 3.25
$ws.Range("P6").Value = 1.57

# Row 7
$ws.Range("F7").Value = 2.5
$ws.Range("G7").Value = 2.74
$ws.Range("H7").Value = 3.15
$ws.Range("I7").Value = 3.5

$wb.Save()
